$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Refresh the cached "datetimeFigureOut" footer field on every slide
#    layout: 6/7/2018 -> 21/10/2018
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shp = $layout.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "6/7/2018") {
                $tr.Text = "21/10/2018"
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Slide 1 content edits (rename "address book" -> "task book" domain
#    terms in the activity diagram)
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# 2a. TextBox 47 -- "[command commits address book]" -> "[command commits task book]"
$shp7 = $s.Shapes.Item(7)
$origHeight7 = $shp7.Height
$tr7 = $shp7.TextFrame.TextRange
$run7b = $tr7.Characters(2, 29)
$run7b.Text = "command commits task book]"
# Editing the run re-triggers the textbox's "shrink/grow to fit" autosize;
# put the height back so only the text actually changes.
$shp7.Height = $origHeight7 + 0.00004

# 2b. Rectangle: Rounded Corners 50 --
#     "Purge redundant states and then save address book to addressBookStateList "
#     -> "Purge redundant states and then save task book to taskBookStateList "
$shp8 = $s.Shapes.Item(8)
$tr8 = $shp8.TextFrame.TextRange

# Edit the right-hand run first (addressBookStateList -> taskBookStateList) so
# the character offsets of the still-untouched left run stay valid.
$run8b = $tr8.Characters(54, 20)
$run8b.Text = "taskBookStateList"

$tr8again = $shp8.TextFrame.TextRange
$run8a = $tr8again.Characters(1, 53)
$run8a.Text = "Purge redundant states and then save task book to "
